$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.502.32'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '2.285.09'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.64'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.43'
$ws.Range('E6').Value = '  +7.01%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +3.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.29'
$ws.Range('E10').Value = '  +12.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.74'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').Value = '2.640.64'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.56'
$ws.Range('E15').Value = '  +2.56%  '
$ws.Range('D16').Value = '2.293.65'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.805'
$ws.Range('E17').Value = '  +5.48%  '
$ws.Range('D18').Value = '42.397.89'
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.70'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.00'
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.98'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.98'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.33'
$ws.Range('E28').Value = '  +11.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.59'
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.13'
$ws.Range('E30').Value = '  +2.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.15'
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.29'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  +4.79%  '
$ws.Range('E35').Value = '  +1.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.34'
$ws.Range('E36').Value = '  +2.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.107'
$ws.Range('E37').Value = '  +2.77%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.38'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').Value = '  +4.19%  '
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('E41').Value = '  +6.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.38'
$ws.Range('E42').Value = '  +13.77%  '
$ws.Range('D43').Value = '2.003.77'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.26'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.04'
$ws.Range('E46').Value = '  +5.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.12'
$ws.Range('E47').Value = '  -2.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.70'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.67'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '93.30'
$ws.Range('E51').Value = '  +2.49%  '
